$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 12)
$ws.Range("A12").Value = "The Unreasonable Effectiveness of Deep Features as a Perceptual Metric"
$ws.Range("B12").Value = "https://arxiv.org/abs/1801.03924?utm_source=chatgpt.com"

# Adjust column B width to fit the new (longer) URL content
$ws.Columns.Item(2).ColumnWidth = 53.166666666666664

# Update selection to match final state
$ws.Range("F14").Select()
